$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cycle_2021-2022")

# Copy formatting (row style / number format) from the last block of
# existing data rows down onto the new rows (2 new 15-row date blocks),
# then fill in the actual values/formulas cell by cell below.
$ws.Range("A89:E103").Copy()
$ws.Range("A104:E118").PasteSpecial(-4122)
$ws.Range("A119:E133").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A104").Value = 44566
$ws.Range("B104").Value = "total applicants"
$ws.Range("C104").Value = 1765
$ws.Range("D104").Value = 3.69
$ws.Range("E104").Value = 511.4

$ws.Range("A105").Value = 44566
$ws.Range("B105").Value = "withdraw before acceptance (WB)"
$ws.Range("C105").Formula = "=NA()"
$ws.Range("D105").Formula = "=NA()"
$ws.Range("E105").Formula = "=NA()"

$ws.Range("A106").Value = 44566
$ws.Range("B106").Value = "rejected groups"
$ws.Range("C106").Formula = "=NA()"
$ws.Range("D106").Formula = "=NA()"
$ws.Range("E106").Formula = "=NA()"

$ws.Range("A107").Value = 44566
$ws.Range("B107").Value = "preliminary rejection"
$ws.Range("C107").Formula = "=NA()"
$ws.Range("D107").Formula = "=NA()"
$ws.Range("E107").Formula = "=NA()"

$ws.Range("A108").Value = 44566
$ws.Range("B108").Value = "passive withdrawal"
$ws.Range("C108").Value = 1
$ws.Range("D108").Value = 4
$ws.Range("E108").Value = 524

$ws.Range("A109").Value = 44566
$ws.Range("B109").Value = "rejection"
$ws.Range("C109").Formula = "=NA()"
$ws.Range("D109").Formula = "=NA()"
$ws.Range("E109").Formula = "=NA()"

$ws.Range("A110").Value = 44566
$ws.Range("B110").Value = "defer to MD app"
$ws.Range("C110").Formula = "=NA()"
$ws.Range("D110").Formula = "=NA()"
$ws.Range("E110").Formula = "=NA()"

$ws.Range("A111").Value = 44566
$ws.Range("B111").Value = "at least 1 MD/PhD acceptance"
$ws.Range("C111").Value = 338
$ws.Range("D111").Value = 517.5
$ws.Range("E111").Value = 3.84

$ws.Range("A112").Value = 44566
$ws.Range("B112").Value = "available active"
$ws.Range("C112").Formula = "=NA()"
$ws.Range("D112").Formula = "=NA()"
$ws.Range("E112").Formula = "=NA()"

$ws.Range("A113").Value = 44566
$ws.Range("B113").Value = "request secondary"
$ws.Range("C113").Formula = "=NA()"
$ws.Range("D113").Formula = "=NA()"
$ws.Range("E113").Formula = "=NA()"

$ws.Range("A114").Value = 44566
$ws.Range("B114").Value = "interview scheduled"
$ws.Range("C114").Formula = "=NA()"
$ws.Range("D114").Formula = "=NA()"
$ws.Range("E114").Formula = "=NA()"

$ws.Range("A115").Value = 44566
$ws.Range("B115").Value = "available passive"
$ws.Range("C115").Formula = "=NA()"
$ws.Range("D115").Formula = "=NA()"
$ws.Range("E115").Formula = "=NA()"

$ws.Range("A116").Value = 44566
$ws.Range("B116").Value = "no action"
$ws.Range("C116").Formula = "=NA()"
$ws.Range("D116").Formula = "=NA()"
$ws.Range("E116").Formula = "=NA()"

$ws.Range("A117").Value = 44566
$ws.Range("B117").Value = "hold"
$ws.Range("C117").Formula = "=NA()"
$ws.Range("D117").Formula = "=NA()"
$ws.Range("E117").Formula = "=NA()"

$ws.Range("A118").Value = 44566
$ws.Range("B118").Value = "available"
$ws.Range("C118").Formula = "=NA()"
$ws.Range("D118").Formula = "=NA()"
$ws.Range("E118").Formula = "=NA()"

$ws.Range("A119").Value = 44576
$ws.Range("B119").Value = "total applicants"
$ws.Range("C119").Value = 1766
$ws.Range("D119").Value = 3.69
$ws.Range("E119").Value = 511.4

$ws.Range("A120").Value = 44576
$ws.Range("B120").Value = "withdraw before acceptance (WB)"
$ws.Range("C120").Value = 40
$ws.Range("D120").Value = 3.73
$ws.Range("E120").Value = 512.5

$ws.Range("A121").Value = 44576
$ws.Range("B121").Value = "rejected groups"
$ws.Range("C121").Value = 1051
$ws.Range("D121").Value = 3.65
$ws.Range("E121").Value = 510

$ws.Range("A122").Value = 44576
$ws.Range("B122").Value = "preliminary rejection"
$ws.Range("C122").Formula = "=NA()"
$ws.Range("D122").Formula = "=NA()"
$ws.Range("E122").Formula = "=NA()"

$ws.Range("A123").Value = 44576
$ws.Range("B123").Value = "passive withdrawal"
$ws.Range("C123").Value = 2
$ws.Range("D123").Value = 3.99
$ws.Range("E123").Value = 521

$ws.Range("A124").Value = 44576
$ws.Range("B124").Value = "rejection"
$ws.Range("C124").Formula = "=NA()"
$ws.Range("D124").Formula = "=NA()"
$ws.Range("E124").Formula = "=NA()"

$ws.Range("A125").Value = 44576
$ws.Range("B125").Value = "defer to MD app"
$ws.Range("C125").Formula = "=NA()"
$ws.Range("D125").Formula = "=NA()"
$ws.Range("E125").Formula = "=NA()"

$ws.Range("A126").Value = 44576
$ws.Range("B126").Value = "at least 1 MD/PhD acceptance"
$ws.Range("C126").Value = 378
$ws.Range("D126").Value = 3.84
$ws.Range("E126").Value = 517.4

$ws.Range("A127").Value = 44576
$ws.Range("B127").Value = "available active"
$ws.Range("C127").Value = 165
$ws.Range("D127").Value = 3.67
$ws.Range("E127").Value = 511.4

$ws.Range("A128").Value = 44576
$ws.Range("B128").Value = "request secondary"
$ws.Range("C128").Formula = "=NA()"
$ws.Range("D128").Formula = "=NA()"
$ws.Range("E128").Formula = "=NA()"

$ws.Range("A129").Value = 44576
$ws.Range("B129").Value = "interview scheduled"
$ws.Range("C129").Formula = "=NA()"
$ws.Range("D129").Formula = "=NA()"
$ws.Range("E129").Formula = "=NA()"

$ws.Range("A130").Value = 44576
$ws.Range("B130").Value = "available passive"
$ws.Range("C130").Value = 132
$ws.Range("D130").Value = 3.58
$ws.Range("E130").Value = 504

$ws.Range("A131").Value = 44576
$ws.Range("B131").Value = "no action"
$ws.Range("C131").Formula = "=NA()"
$ws.Range("D131").Formula = "=NA()"
$ws.Range("E131").Formula = "=NA()"

$ws.Range("A132").Value = 44576
$ws.Range("B132").Value = "hold"
$ws.Range("C132").Formula = "=NA()"
$ws.Range("D132").Formula = "=NA()"
$ws.Range("E132").Formula = "=NA()"

$ws.Range("A133").Value = 44576
$ws.Range("B133").Value = "available"
$ws.Range("C133").Formula = "=NA()"
$ws.Range("D133").Formula = "=NA()"
$ws.Range("E133").Formula = "=NA()"

# Scroll the frozen-header view down to the new bottom rows and select
# the cell that was active there after entering the new data.
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("H132").Select()
